# Quarterly indexing esoteric bug-fix operation
#
# Column A holds the "as-of" date (Excel serial date number, formatted as a
# date by style s="1") for each forecast row. These dates were erroneously
# stamped with the first day of the reference quarter's first month
# (e.g. 2010-04-01). The correct "as-of" date is the 15th of the month
# following that (e.g. 2010-05-15), which is when the underlying source
# data actually becomes available. Re-stamp every row in column A (rows
# 2 through 63) with the corrected serial date value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFixes = @{
    2  = 25614
    3  = 25614
    4  = 25614
    5  = 25614
    6  = 25614
    7  = 25614
    8  = 25614
    9  = 25614
    10 = 25614
    11 = 25614
    12 = 40313
    13 = 40405
    14 = 40497
    15 = 40589
    16 = 40678
    17 = 40770
    18 = 40862
    19 = 40954
    20 = 41044
    21 = 41136
    22 = 41228
    23 = 41320
    24 = 41409
    25 = 41501
    26 = 41593
    27 = 41685
    28 = 41774
    29 = 41866
    30 = 41958
    31 = 42050
    32 = 42139
    33 = 42231
    34 = 42323
    35 = 42415
    36 = 42505
    37 = 42597
    38 = 42689
    39 = 42781
    40 = 42870
    41 = 42962
    42 = 43054
    43 = 43146
    44 = 43235
    45 = 43327
    46 = 43419
    47 = 43511
    48 = 43600
    49 = 43692
    50 = 43784
    51 = 43876
    52 = 43966
    53 = 44058
    54 = 44150
    55 = 44242
    56 = 44331
    57 = 44423
    58 = 44515
    59 = 44607
    60 = 44696
    61 = 44788
    62 = 44880
    63 = 44972
}

foreach ($row in $dateFixes.Keys) {
    $ws.Cells.Item($row, 1).Value = $dateFixes[$row]
}
